# Applies the data update described by the commit "add almost all lineal".
#
# Only numeric-looking *text* values change (the source workbook stores
# every expression/number on these sheets as literal text, not as real
# numbers) on the sheets:
#   - Restricciones_del_follower
#   - Punto_modificado
#   - Vector_bf
#   - Vector_BF
#   - Vector_Alpha (this one's last two rows are genuine numbers)
#
# Because typing a numeric-looking string into a cell normally makes Excel
# coerce it to a real number, each text cell is temporarily switched to the
# "Text" number format before the value is assigned, then restored to the
# "Normal" style afterwards so no visible formatting change is left behind.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2 (J_0_L0_v)
Set-TextValue $ws "A2" "0.9270360882359463y_1 + 0.989582907613092y_2"
Set-TextValue $ws "B2" "6.4258418245948965"
Set-TextValue $ws "D2" "0.0866877650392671"
Set-TextValue $ws "E2" "0.18254418018572394"
Set-TextValue $ws "F2" "0.7018207991573857"

# Row 3 (J_0_L0_v)
Set-TextValue $ws "A3" "-4 + 0.8965435170293287y_1 - 0.05312758170671375y_2"
Set-TextValue $ws "B3" "-0.19754896467942032"
Set-TextValue $ws "D3" "0.9648587319705634"
Set-TextValue $ws "E3" "0"
Set-TextValue $ws "F3" "0.425069827518858"

# Row 4 (J_0_LP_v)
Set-TextValue $ws "A4" "-16 - 2x - 0.32852662910738983y_1 + 3.3177671905060135y_2"
Set-TextValue $ws "B4" "-21.269467588355326"
Set-TextValue $ws "D4" "0.9761226555169311"
Set-TextValue $ws "E4" "0.7333987827273575"
Set-TextValue $ws "F4" "0"

# Row 5 (J_Ne_L0_v)
Set-TextValue $ws "A5" "-48 + 8x + 0.20867847500363246y_1 - 0.4063640843797479y_2"
Set-TextValue $ws "B5" "-1.0489988692656744"
Set-TextValue $ws "D5" "0.8143958706897286"
Set-TextValue $ws "E5" "0.3400618792126797"
Set-TextValue $ws "F5" "0"

# Row 6 (J_Ne_L0_v)
Set-TextValue $ws "A6" "12 - 2x - 0.4173569500072649y_1 + 0.8127281687594958y_2"
Set-TextValue $ws "B6" "0.35976267716903676"
Set-TextValue $ws "D6" "0.23927405565041526"
Set-TextValue $ws "E6" "0.8897524486363444"
Set-TextValue $ws "F6" "0"

# ---------------------------------------------------------------
# Punto_modificado
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws "A2" "5.875840352759835"
Set-TextValue $ws "B2" "4.382729079133727"
Set-TextValue $ws "C2" "2.387767396848251"

# ---------------------------------------------------------------
# Vector_bf (sheet index 5). Using the numeric index rather than the
# sheet name because "Vector_bf" / "Vector_BF" only differ by case and
# name-based lookup in this runtime is case-insensitive (it would
# otherwise resolve both names to the same sheet).
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws "A2" "0.30519755975240304"
Set-TextValue $ws "A3" "-3.1365953726124545"

# ---------------------------------------------------------------
# Vector_BF (sheet index 6)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws "A2" "1.525807429025966"
Set-TextValue $ws "A3" "3.3720967610146593"
Set-TextValue $ws "A4" "-5.198826963536688"

# ---------------------------------------------------------------
# Vector_Alpha -- these two cells are genuine numbers (not text) in
# the workbook, so assign them directly as numeric values.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 0.4890013063486953
$ws.Range("A3").Value = 0.2511148273336813
